# Update resources/ResourceFile_Method_HT.xlsx:
#   - prevalence2018!C21:C122 revised in line with STEP data
#   - incidence2018_plus!C21:C122 are =prevalence2018!Cxx/100 and recalc automatically
#   - sheet view / selection state updated to match author's last position
#   - "data" sheet selection updated

$wb = $excel.ActiveWorkbook

# --- prevalence2018 -------------------------------------------------------
$ws = $wb.Worksheets.Item("prevalence2018")
$ws.Activate()

# Ages 19-24 (rows 21-26): 0.35 -> 0
$ws.Range("C21:C26").Value = 0

# Ages 25-34 (rows 27-36): 0.35 -> 0.36
$ws.Range("C27:C36").Value = 0.36

# Ages 35-44 (rows 37-46): 0.43 -> 0.425
$ws.Range("C37:C46").Value = 0.425

# Ages 45-54 (rows 47-56): 0.57 -> 0.5
$ws.Range("C47:C56").Value = 0.5

# Ages 55+ (rows 57-122): 0.9 -> 0.97
$ws.Range("C57:C122").Value = 0.97

# Sheet view: scrolled up a bit, new selection block
$ws.Range("C57:C122").Select()

# --- data sheet -------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")
$wsData.Activate()
$wsData.Range("D5:D8").Select()

# Restore prevalence2018 as the active/visible sheet (matches activeTab in workbook.xml)
$ws.Activate()
